$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds dates that need to move from 2023-09-16 (serial 45185)
# to 2023-10-05 (serial 45204) for every data row (rows 2 through 83).
for ($row = 2; $row -le 83; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45185) {
        $cell.Value2 = 45204
    }
}
